$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.820.66'
$ws.Range("E2").Value = '  -2.20%  '

$ws.Range("D3").Value = '1.620.81'
$ws.Range("E3").Value = '  -1.83%  '

$ws.Range("E4").Value = '  -0.52%  '

$ws.Range("D5").Value = '308.28'
$ws.Range("E5").Value = '  -1.26%  '

$ws.Range("E6").Value = '  -0.41%  '

$ws.Range("D7").Value = '0.3929'
$ws.Range("E7").Value = '  +0.35%  '

$ws.Range("D8").Value = '0.3849'
$ws.Range("E8").Value = '  -1.64%  '

$ws.Range("D9").Value = '1.001'
$ws.Range("E9").Value = '  -0.51%  '

$ws.Range("D10").Value = '49.62'
$ws.Range("E10").Value = '  -2.00%  '

$ws.Range("D11").Value = '1.353'
$ws.Range("E11").Value = '  -1.79%  '

$ws.Range("D12").Value = '0.08469'
$ws.Range("E12").Value = '  -0.96%  '

$ws.Range("D13").Value = '23.72'
$ws.Range("E13").Value = '  -4.75%  '

$ws.Range("D14").Value = '7.065'
$ws.Range("E14").Value = '  -1.78%  '

$ws.Range("D15").Value = '7.622'
$ws.Range("E15").Value = '  +0.24%  '

$ws.Range("E16").Value = '  -1.29%  '

$ws.Range("D17").Value = '1.614.21'
$ws.Range("E17").Value = '  -2.97%  '

$ws.Range("D18").Value = '93.87'
$ws.Range("E18").Value = '  +0.77%  '

$ws.Range("D19").Value = '0.06931'
$ws.Range("E19").Value = '  -0.40%  '

$ws.Range("D20").Value = '19.98'
$ws.Range("E20").Value = '  -4.86%  '

$ws.Range("E21").Value = '  -2.47%  '

$ws.Range("D22").Value = '0.9995'
$ws.Range("E22").Value = '  -0.59%  '

$ws.Range("E23").Value = '  -2.40%  '

$ws.Range("D24").Value = '23.829.44'
$ws.Range("E24").Value = '  -2.15%  '

$ws.Range("D25").Value = '2.488'
$ws.Range("E25").Value = '  +5.72%  '

$ws.Range("D26").Value = '2.834'
$ws.Range("E26").Value = '  +2.60%  '

$ws.Range("D27").Value = '22.24'
$ws.Range("E27").Value = '  -1.85%  '

$ws.Range("D28").Value = '156.93'
$ws.Range("E28").Value = '  -1.03%  '

$ws.Range("D29").Value = '140.45'
$ws.Range("E29").Value = '  -3.13%  '

$ws.Range("D30").Value = '5.303'
$ws.Range("E30").Value = '  -7.81%  '

$ws.Range("D31").Value = '7.807'
$ws.Range("E31").Value = '  -3.43%  '

$ws.Range("D32").Value = '2.478'
$ws.Range("E32").Value = '  -1.20%  '

$ws.Range("D33").Value = '1.791.01'
$ws.Range("E33").Value = '  -2.15%  '

$ws.Range("D34").Value = '0.08135'
$ws.Range("E34").Value = '  -0.43%  '

$ws.Range("D35").Value = '0.9903'
$ws.Range("E35").Value = '  -1.94%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02898'
$ws.Range("E36").Value = '  -3.99%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '6.632'
$ws.Range("E37").Value = '  -3.01%  '

$ws.Range("D38").Value = '0.2670'
$ws.Range("E38").Value = '  -3.28%  '

$ws.Range("D39").Value = '0.09162'
$ws.Range("E39").Value = '  -3.97%  '

$ws.Range("D40").Value = '10.33'
$ws.Range("E40").Value = '  +1.57%  '

$ws.Range("D41").Value = '13.67'
$ws.Range("E41").Value = '  +2.90%  '

$ws.Range("D42").Value = '1.425'
$ws.Range("E42").Value = '  -4.44%  '

$ws.Range("D43").Value = '0.7532'
$ws.Range("E43").Value = '  -2.81%  '

$ws.Range("D44").Value = '16.01'
$ws.Range("E44").Value = '  -1.32%  '

$ws.Range("D45").Value = '0.6928'
$ws.Range("E45").Value = '  -0.94%  '

$ws.Range("D46").Value = '2.474'
$ws.Range("E46").Value = '  -3.03%  '

$ws.Range("D47").Value = '4.073'
$ws.Range("E47").Value = '  -1.89%  '

$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  -0.48%  '

$ws.Range("D49").Value = '0.08245'
$ws.Range("E49").Value = '  -3.66%  '

$ws.Range("D50").Value = '135.66'
$ws.Range("E50").Value = '  -0.60%  '

$ws.Range("D51").Value = '1.199'
$ws.Range("E51").Value = '  -7.60%  '
